$d = $word.ActiveDocument

# --- Part 1: header block (first four paragraphs) ---
# Remove right alignment and change font color from yellow (FFFF00) to orange (FF6600)
for ($i = 1; $i -le 4; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Alignment = 0
    $p.Range.Font.Color = 26367
}

# --- Part 2: move the _GoBack bookmark from the middle of "riverbank" to
#     right after "A " in the "A Cat, a Parrot, and a Bag of Seed:" heading,
#     and merge the now-needlessly-split "riverba" / "nk..." runs back
#     together, while leaving the following "s" / ". Show how..." runs
#     (which are unrelated to the bookmark move) untouched. ---

# Locate a stable anchor straddling the run boundaries that must stay put so
# we can compute their character offsets without hard-coding the (long,
# non-breaking-space-laden) paragraph text.
$anchor = $d.Content
$anchor.Find.Execute("bag of seeds. Show", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$mergeBoundary = $anchor.Start + ("bag of seed").Length
$sRunEnd = $mergeBoundary + 1

# Temporary "barrier" bookmarks keep the edit below from coalescing the
# unrelated "s" run and the ". Show how..." run into the merged sentence.
$d.Bookmarks.Add("ZZZ_barrier1", $d.Range($mergeBoundary, $mergeBoundary))
$d.Bookmarks.Add("ZZZ_barrier2", $d.Range($sRunEnd, $sRunEnd))

# Remove the bookmark from its old spot in the middle of "riverbank".
$d.Bookmarks.Item("_GoBack").Delete()

# Touch the text spanning the old run split so the document engine
# recombines "riverba" + "nk..." into a single run again.
$para7Start = $d.Paragraphs.Item(7).Range.Start
$touch = $d.Range($para7Start, $para7Start + 2)
$touchText = $touch.Text
$touch.Delete()
$d.Range($para7Start, $para7Start).InsertAfter($touchText)

# Remove the temporary barrier bookmarks (plain bookmark deletion does not
# itself trigger a run coalesce).
$d.Bookmarks.Item("ZZZ_barrier1").Delete()
$d.Bookmarks.Item("ZZZ_barrier2").Delete()

# Re-add the _GoBack bookmark right after "A " in the heading paragraph.
$bmPos = $d.Paragraphs.Item(6).Range.Start + 2
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

Write-Host "Done"
